$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh: the price-record rows for 2021-06-03 and 2021-03-18 had
# swapped into the wrong rows; this swaps the affected fields (date, volume,
# prices, origin) back between row 2 and row 4. Row 3 is untouched.

$cols = @("D", "J", "K", "L", "M", "O", "P")

$row2Values = @{}
$row4Values = @{}

foreach ($col in $cols) {
    $row2Values[$col] = $ws.Range("$col`2").Value2()
    $row4Values[$col] = $ws.Range("$col`4").Value2()
}

foreach ($col in $cols) {
    $ws.Range("$col`2").Value = $row4Values[$col]
    $ws.Range("$col`4").Value = $row2Values[$col]
}
